$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2195121951219512
$ws.Range("C2").Value = 0.5257452574525745
$ws.Range("J2").Value = 0.01897018970189702
$ws.Range("P2").Value = 0.1626016260162602
$ws.Range("S2").Value = 0.07317073170731707
$ws.Range("C3").Value = 0.02857142857142857
$ws.Range("J3").Value = 0.02857142857142857
$ws.Range("P3").Value = 0.7476190476190476
$ws.Range("S3").Value = 0.1952380952380952
$ws.Range("J4").Value = 0.03846153846153846
$ws.Range("P4").Value = 0.6538461538461539
$ws.Range("S4").Value = 0.3076923076923077
$ws.Range("P5").Value = 1
$ws.Range("B6").Value = 0.0514018691588785
$ws.Range("D6").Value = 0.01869158878504673
$ws.Range("F6").Value = 0.04672897196261682
$ws.Range("J6").Value = 0.2102803738317757
$ws.Range("O6").Value = 0.01401869158878505
$ws.Range("Q6").Value = 0.1962616822429906
$ws.Range("R6").Value = 0.08411214953271028
$ws.Range("S6").Value = 0.3785046728971962
$ws.Range("B7").Value = 0.1586538461538461
$ws.Range("D7").Value = 0.03846153846153846
$ws.Range("F7").Value = 0.0576923076923077
$ws.Range("J7").Value = 0.125
$ws.Range("O7").Value = 0.01923076923076923
$ws.Range("Q7").Value = 0.1682692307692308
$ws.Range("R7").Value = 0.04326923076923077
$ws.Range("S7").Value = 0.3894230769230769
$ws.Range("B8").Value = 0.0928030303030303
$ws.Range("D8").Value = 0.01515151515151515
$ws.Range("F8").Value = 0.07007575757575757
$ws.Range("J8").Value = 0.115530303030303
$ws.Range("O8").Value = 0.02083333333333333
$ws.Range("Q8").Value = 0.1666666666666667
$ws.Range("R8").Value = 0.0928030303030303
$ws.Range("S8").Value = 0.4261363636363636
$ws.Range("B9").Value = 0.08018867924528301
$ws.Range("D9").Value = 0.01886792452830189
$ws.Range("F9").Value = 0.0660377358490566
$ws.Range("J9").Value = 0.1226415094339623
$ws.Range("O9").Value = 0.02358490566037736
$ws.Range("Q9").Value = 0.2028301886792453
$ws.Range("R9").Value = 0.08018867924528301
$ws.Range("S9").Value = 0.4056603773584906
$ws.Range("B10").Value = 0.1186079545454545
$ws.Range("D10").Value = 0.01775568181818182
$ws.Range("E10").Value = 0.0007102272727272727
$ws.Range("F10").Value = 0.0546875
$ws.Range("J10").Value = 0.1392045454545454
$ws.Range("O10").Value = 0.01775568181818182
$ws.Range("Q10").Value = 0.2258522727272727
$ws.Range("R10").Value = 0.07244318181818182
$ws.Range("S10").Value = 0.3529829545454545
$ws.Range("G11").Value = 0.1551724137931035
$ws.Range("J11").Value = 0.05862068965517241
$ws.Range("K11").Value = 0.1724137931034483
$ws.Range("L11").Value = 0.6068965517241379
$ws.Range("S11").Value = 0.006896551724137931
$ws.Range("G12").Value = 0.7653631284916201
$ws.Range("J12").Value = 0.1955307262569832
$ws.Range("K12").Value = 0.0111731843575419
$ws.Range("L12").Value = 0.0223463687150838
$ws.Range("S12").Value = 0.00558659217877095
$ws.Range("G13").Value = 0.6875
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.0625
$ws.Range("F15").Value = 0.02232142857142857
$ws.Range("H15").Value = 0.1607142857142857
$ws.Range("I15").Value = 0.05803571428571429
$ws.Range("J15").Value = 0.3705357142857143
$ws.Range("K15").Value = 0.05803571428571429
$ws.Range("M15").Value = 0.01785714285714286
$ws.Range("O15").Value = 0.08035714285714286
$ws.Range("S15").Value = 0.2321428571428572
$ws.Range("F16").Value = 0.008474576271186441
$ws.Range("H16").Value = 0.211864406779661
$ws.Range("I16").Value = 0.08050847457627118
$ws.Range("J16").Value = 0.4279661016949153
$ws.Range("K16").Value = 0.08050847457627118
$ws.Range("M16").Value = 0.01271186440677966
$ws.Range("O16").Value = 0.02966101694915254
$ws.Range("S16").Value = 0.1483050847457627
$ws.Range("F17").Value = 0.007662835249042145
$ws.Range("H17").Value = 0.1800766283524904
$ws.Range("I17").Value = 0.1053639846743295
$ws.Range("J17").Value = 0.417624521072797
$ws.Range("K17").Value = 0.07279693486590039
$ws.Range("M17").Value = 0.01724137931034483
$ws.Range("O17").Value = 0.06896551724137931
$ws.Range("S17").Value = 0.1302681992337165
$ws.Range("F18").Value = 0.01
$ws.Range("H18").Value = 0.235
$ws.Range("I18").Value = 0.055
$ws.Range("J18").Value = 0.4
$ws.Range("K18").Value = 0.095
$ws.Range("M18").Value = 0.025
$ws.Range("O18").Value = 0.05
$ws.Range("S18").Value = 0.13
$ws.Range("F19").Value = 0.01520637219406227
$ws.Range("H19").Value = 0.220854453294714
$ws.Range("I19").Value = 0.08254887762490949
$ws.Range("J19").Value = 0.3758146270818248
$ws.Range("K19").Value = 0.11151339608979
$ws.Range("M19").Value = 0.01955104996379435
$ws.Range("O19").Value = 0.05792903692976104
$ws.Range("S19").Value = 0.1165821868211441
